$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.722.31'
$ws.Range('E2').Value = '  -0.60%  '
$ws.Range('D3').Value = '1.848.60'
$ws.Range('E3').Value = '  -1.15%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.014'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -2.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '319.41'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.53%  '
$ws.Range('E6').Value = '  -1.93%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4314'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.67%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3751'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.90%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07351'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.70%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8774'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.57%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.63'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.81%  '
$ws.Range('D12').Value = '1.851.89'
$ws.Range('E12').Value = '  -1.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.730'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.90%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.449'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.74%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.07129'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.23%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '89.12'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.85%  '
$ws.Range('E17').Value = '  -2.11%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009005'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.61%  '
$ws.Range('E19').Value = '  -2.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.49'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.91%  '
$ws.Range('D21').Value = '27.721.81'
$ws.Range('E21').Value = '  -0.72%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.220'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.10'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.29%  '
$ws.Range('D24').Value = '2.078.78'
$ws.Range('E24').Value = '  -1.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.002'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.11%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '155.33'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.67'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.35%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.175'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +9.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.393'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.33%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '119.31'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.67%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08947'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.42%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.234'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.51%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7780'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.69%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.563'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.21%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.911'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.34%  '
$ws.Range('E36').Value = '  -1.96%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.134'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05348'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.58%  '
$ws.Range('E39').Value = '  -0.80%  '
$ws.Range('E40').Value = '  +5.64%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.891'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5141'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.80%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1689'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.70%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '8.811'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.31%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.77'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.04%  '
$ws.Range('B46').Value = 'Quant'
$ws.Range('C46').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '108.94'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.86%  '
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.4775'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.44%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.06476'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.20%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.694'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.85%  '
$ws.Range('B50').Value = 'PaxDollar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.012'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.02%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.857'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.13%  '
